$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") bumped by one day (46059 -> 46060) for every data row (2-11)
foreach ($r in 2..11) {
    $ws.Cells.Item($r, 3).Value = 46060
}

# Rows 7-11 were re-ordered (their Beteckning/Datum/Area values rotated),
# while other columns on those rows stayed the same.
$ws.Range("A7").Value = "A 33037-2025"
$ws.Range("B7").Value = 45840.39623842593
$ws.Range("G7").Value = 0.8

$ws.Range("A8").Value = "A 33033-2025"
$ws.Range("B8").Value = 45840.39188657407
$ws.Range("G8").Value = 0.7

$ws.Range("A9").Value = "A 6314-2022"
$ws.Range("B9").Value = 44600
$ws.Range("G9").Value = 3

$ws.Range("A10").Value = "A 25610-2024"
$ws.Range("B10").Value = 45463
$ws.Range("G10").Value = 2.9

$ws.Range("A11").Value = "A 57810-2022"
$ws.Range("B11").Value = 44897
$ws.Range("G11").Value = 3.3
